$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 90, pushing existing rows 90:217 down to 91:218.
$ws.Rows("90:90").Insert()

# Populate the newly inserted row 90 with the new data record.
$ws.Range("A90").Value = 3
$ws.Range("B90").Value = "Femacal de La Calera"
$ws.Range("C90").Value = "Coquimbo"
$ws.Range("D90").Value = 44579
$ws.Range("E90").Value = 5
$ws.Range("F90").Value = 100112001
$ws.Range("G90").Value = "Berenjena"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 55
$ws.Range("K90").Value = 10000
$ws.Range("L90").Value = 10000
$ws.Range("M90").Value = 10000
$ws.Range("N90").Value = "$/caja 60 unidades"
$ws.Range("O90").Value = "Región de Arica y Parinacota"
$ws.Range("P90").Value = 167
$ws.Range("Q90").Value = 60
$ws.Range("R90").Value = "Hortaliza"
